$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns M (2021) and N (2022) ------------------------------------
# Copy the formatting of the existing "2020" column (L) onto the two new
# columns before writing values, so the new cells pick up the same number
# format / styles the rest of the table uses.
$ws.Range("L4:L12").Copy()
$ws.Range("M4:N12").PasteSpecial(-4122)  # xlPasteFormats

# Year header row
$ws.Range("M4").Value = 2021
$ws.Range("N4").Value = 2022

# Data rows (2021 repeats the 2020 figure, 2022 is the new figure)
$ws.Range("M5").Value = 5.6
$ws.Range("N5").Value = 6.3

$ws.Range("M6").Value = 0.8
$ws.Range("N6").Value = 0.8

$ws.Range("M7").Value = 1.9
$ws.Range("N7").Value = 2.4

$ws.Range("M8").Value = 0.7
$ws.Range("N8").Value = 0.7

$ws.Range("M9").Value = 0.7
$ws.Range("N9").Value = 0.8

$ws.Range("M10").Value = 0.9
$ws.Range("N10").Value = 1

$ws.Range("M11").Value = 0.3
$ws.Range("N11").Value = 0.2

$ws.Range("M12").Value = 0.2
$ws.Range("N12").Value = 0.4

# --- New footnote row 14 ----------------------------------------------------
# Row 13 already carries a wrapped footnote style (B13/C13, style id 17).
# Copy that formatting down to the new row and add the new footnote text.
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B14").Value = "По данным лесоустройства 2022 года Лесной службы при Министерстве чрезвычайных ситуаций КР"
$ws.Rows.Item(14).RowHeight = 34.5
